$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "66.887.92"
Set-TextValue "E2" "  +2.23%  "
Set-TextValue "D3" "3.097.36"
Set-TextValue "E3" "  +4.84%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "579.54"
Set-TextValue "E5" "  +1.50%  "
Set-TextValue "D6" "173.30"
Set-TextValue "E6" "  +7.71%  "
Set-TextValue "E7" "  -0.03%  "
Set-TextValue "D8" "3.092.42"
Set-TextValue "E8" "  +4.73%  "
Set-TextValue "E9" "  +1.17%  "
Set-TextValue "E10" "  -2.66%  "
Set-TextValue "D11" "0.156"
Set-TextValue "E11" "  +3.74%  "
Set-TextValue "E12" "  +4.86%  "
Set-TextValue "E13" "  +2.15%  "
Set-TextValue "D14" "37.14"
Set-TextValue "E14" "  +7.37%  "
Set-TextValue "E15" "  -0.21%  "
Set-TextValue "D16" "3.608.90"
Set-TextValue "E16" "  +4.86%  "
Set-TextValue "D17" "66.853.59"
Set-TextValue "E17" "  +2.27%  "
Set-TextValue "D18" "7.19"
Set-TextValue "E18" "  +2.06%  "
Set-TextValue "D19" "3.097.86"
Set-TextValue "E19" "  +5.02%  "
Set-TextValue "D20" "16.20"
Set-TextValue "E20" "  +1.15%  "
Set-TextValue "D21" "482.72"
Set-TextValue "E21" "  +8.31%  "
Set-TextValue "E22" "  +2.32%  "
Set-TextValue "D23" "7.51"
Set-TextValue "E23" "  +2.69%  "
Set-TextValue "D24" "84.08"
Set-TextValue "E24" "  +1.98%  "
Set-TextValue "E25" "  +4.36%  "
Set-TextValue "E26" "  +6.27%  "
Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  -0.10%  "
Set-TextValue "D28" "10.00"
Set-TextValue "E28" "  -0.03%  "
Set-TextValue "D29" "7.96"
Set-TextValue "E29" "  +0.01%  "
Set-TextValue "D30" "2.39"
Set-TextValue "E30" "  -3.70%  "
Set-TextValue "E31" "  +3.27%  "
Set-TextValue "E32" "  -0.84%  "
Set-TextValue "D33" "28.81"
Set-TextValue "E33" "  +5.87%  "
Set-TextValue "E34" "  +1.79%  "
Set-TextValue "D35" "0.999"
Set-TextValue "E35" "  +0.11%  "
Set-TextValue "D36" "1.00"
Set-TextValue "E36" "  +2.95%  "
Set-TextValue "D37" "5.88"
Set-TextValue "E37" "  +2.35%  "
Set-TextValue "D38" "47.92"
Set-TextValue "E38" "  +7.04%  "
Set-TextValue "D39" "2.12"
Set-TextValue "E39" "  +7.22%  "
Set-TextValue "D40" "50.17"
Set-TextValue "E40" "  +2.06%  "
Set-TextValue "D41" "0.316"
Set-TextValue "E41" "  +4.90%  "
Set-TextValue "E42" "  +0.53%  "
Set-TextValue "E43" "  +1.15%  "
Set-TextValue "E44" "  -1.41%  "
Set-TextValue "E45" "  +2.16%  "
Set-TextValue "D46" "2.816.11"
Set-TextValue "E46" "  +4.94%  "
Set-TextValue "D47" "378.70"
Set-TextValue "E47" "  -1.22%  "
Set-TextValue "D48" "134.98"
Set-TextValue "E48" "  +1.14%  "
Set-TextValue "E49" "  -0.01%  "
Set-TextValue "D50" "24.79"
Set-TextValue "E50" "  +4.77%  "
Set-TextValue "E51" "  +1.63%  "
